$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "BTC"
$ws.Range("C2").Value = "Bitcoin"
$ws.Range("D2").Value = 44074
$ws.Range("E2").Value = 864386053192
$ws.Range("F2").Value = 11043791530
$ws.Range("G2").Value = 0.87121

$ws.Range("B3").Value = "ETH"
$ws.Range("C3").Value = "Ethereum"
$ws.Range("D3").Value = 2237.69
$ws.Range("E3").Value = 269212062847
$ws.Range("F3").Value = 7942485475
$ws.Range("G3").Value = -0.04676

$ws.Range("B4").Value = "USDT"
$ws.Range("C4").Value = "Tether"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 93690138972
$ws.Range("F4").Value = 25149150612
$ws.Range("G4").Value = -0.02817

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "BNB"
$ws.Range("D5").Value = 306.08
$ws.Range("E5").Value = 47156237475
$ws.Range("F5").Value = 590373742
$ws.Range("G5").Value = 0.07722999999999999

$ws.Range("B6").Value = "SOL"
$ws.Range("C6").Value = "Solana"
$ws.Range("D6").Value = 94.37
$ws.Range("E6").Value = 40725519860
$ws.Range("F6").Value = 2034755507
$ws.Range("G6").Value = -0.30898

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "XRP"
$ws.Range("D7").Value = 0.565615
$ws.Range("E7").Value = 30690989323
$ws.Range("F7").Value = 684380370
$ws.Range("G7").Value = -0.32723

$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "USDC"
$ws.Range("D8").Value = 0.999676
$ws.Range("E8").Value = 25400638532
$ws.Range("F8").Value = 4009967635
$ws.Range("G8").Value = 0.009860000000000001

$ws.Range("B9").Value = "STETH"
$ws.Range("C9").Value = "Lido Staked Ether"
$ws.Range("D9").Value = 2235.63
$ws.Range("E9").Value = 20696735176
$ws.Range("F9").Value = 10092396
$ws.Range("G9").Value = -0.15267

$ws.Range("B10").Value = "ADA"
$ws.Range("C10").Value = "Cardano"
$ws.Range("D10").Value = 0.520087
$ws.Range("E10").Value = 18247800940
$ws.Range("F10").Value = 383833565
$ws.Range("G10").Value = -0.32159

$ws.Range("B11").Value = "AVAX"
$ws.Range("C11").Value = "Avalanche"
$ws.Range("D11").Value = 34.77
$ws.Range("E11").Value = 12768889978
$ws.Range("F11").Value = 599429667
$ws.Range("G11").Value = 0.21768

$ws.Range("B12").Value = "DOGE"
$ws.Range("C12").Value = "Dogecoin"
$ws.Range("D12").Value = 0.079929
$ws.Range("E12").Value = 11405914852
$ws.Range("F12").Value = 406016514
$ws.Range("G12").Value = -1.07632

$ws.Range("B13").Value = "DOT"
$ws.Range("C13").Value = "Polkadot"
$ws.Range("D13").Value = 7.18
$ws.Range("E13").Value = 9446806079
$ws.Range("F13").Value = 221049552
$ws.Range("G13").Value = -0.14984

$ws.Range("B14").Value = "TRX"
$ws.Range("C14").Value = "TRON"
$ws.Range("D14").Value = 0.10359
$ws.Range("E14").Value = 9149646600
$ws.Range("F14").Value = 241935838
$ws.Range("G14").Value = -0.09933

$ws.Range("B15").Value = "MATIC"
$ws.Range("C15").Value = "Polygon"
$ws.Range("D15").Value = 0.830101
$ws.Range("E15").Value = 7712435868
$ws.Range("F15").Value = 309242215
$ws.Range("G15").Value = 0.98666

$ws.Range("B16").Value = "LINK"
$ws.Range("C16").Value = "Chainlink"
$ws.Range("D16").Value = 13.5
$ws.Range("E16").Value = 7673317699
$ws.Range("F16").Value = 340513217
$ws.Range("G16").Value = -0.34877

$ws.Range("B17").Value = "TON"
$ws.Range("C17").Value = "Toncoin"
$ws.Range("D17").Value = 2.21
$ws.Range("E17").Value = 7636013413
$ws.Range("F17").Value = 48563694
$ws.Range("G17").Value = 3.40942

$ws.Range("B18").Value = "WBTC"
$ws.Range("C18").Value = "Wrapped Bitcoin"
$ws.Range("D18").Value = 44133
$ws.Range("E18").Value = 6984550433
$ws.Range("F18").Value = 139385498
$ws.Range("G18").Value = 0.89237

$ws.Range("B19").Value = "SHIB"
$ws.Range("C19").Value = "Shiba Inu"
$ws.Range("D19").Value = 0.00000947
$ws.Range("E19").Value = 5571906937
$ws.Range("F19").Value = 100097183
$ws.Range("G19").Value = -1.13943

$ws.Range("B20").Value = "ICP"
$ws.Range("C20").Value = "Internet Computer"
$ws.Range("D20").Value = 11.89
$ws.Range("E20").Value = 5419149196
$ws.Range("F20").Value = 149874381
$ws.Range("G20").Value = -1.44022

$ws.Range("B21").Value = "DAI"
$ws.Range("C21").Value = "Dai"
$ws.Range("D21").Value = 0.996938
$ws.Range("E21").Value = 5256626275
$ws.Range("F21").Value = 213913316
$ws.Range("G21").Value = -0.22965

$ws.Range("B22").Value = "LTC"
$ws.Range("C22").Value = "Litecoin"
$ws.Range("D22").Value = 65.40000000000001
$ws.Range("E22").Value = 4838776321
$ws.Range("F22").Value = 386254645
$ws.Range("G22").Value = 0.56813

$ws.Range("B23").Value = "UNI"
$ws.Range("C23").Value = "Uniswap"
$ws.Range("D23").Value = 6.3
$ws.Range("E23").Value = 4749199001
$ws.Range("F23").Value = 204065391
$ws.Range("G23").Value = 0.55937

$ws.Range("B24").Value = "BCH"
$ws.Range("C24").Value = "Bitcoin Cash"
$ws.Range("D24").Value = 235.53
$ws.Range("E24").Value = 4630937966
$ws.Range("F24").Value = 135763306
$ws.Range("G24").Value = -0.12728

$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "LEO Token"
$ws.Range("D25").Value = 4.03
$ws.Range("E25").Value = 3738069445
$ws.Range("F25").Value = 1013178
$ws.Range("G25").Value = -0.82098

$ws.Range("B26").Value = "ATOM"
$ws.Range("C26").Value = "Cosmos Hub"
$ws.Range("D26").Value = 9.789999999999999
$ws.Range("E26").Value = 3737316528
$ws.Range("F26").Value = 164513430
$ws.Range("G26").Value = -1.13326

$ws.Range("B27").Value = "NEAR"
$ws.Range("C27").Value = "NEAR Protocol"
$ws.Range("D27").Value = 3.36
$ws.Range("E27").Value = 3411741844
$ws.Range("F27").Value = 245125512
$ws.Range("G27").Value = 1.09563

$ws.Range("B28").Value = "XLM"
$ws.Range("C28").Value = "Stellar"
$ws.Range("D28").Value = 0.118299
$ws.Range("E28").Value = 3347630314
$ws.Range("F28").Value = 73161657
$ws.Range("G28").Value = -0.48379

$ws.Range("B29").Value = "OKB"
$ws.Range("C29").Value = "OKB"
$ws.Range("D29").Value = 54.66
$ws.Range("E29").Value = 3285621801
$ws.Range("F29").Value = 3616196
$ws.Range("G29").Value = 0.0766

$ws.Range("B30").Value = "INJ"
$ws.Range("C30").Value = "Injective"
$ws.Range("D30").Value = 37.52
$ws.Range("E30").Value = 3166757556
$ws.Range("F30").Value = 306581159
$ws.Range("G30").Value = -0.35053

$ws.Range("B31").Value = "OP"
$ws.Range("C31").Value = "Optimism"
$ws.Range("D31").Value = 3.37
$ws.Range("E31").Value = 3064274812
$ws.Range("F31").Value = 212211547
$ws.Range("G31").Value = 3.77319

$ws.Range("B32").Value = "FIL"
$ws.Range("C32").Value = "Filecoin"
$ws.Range("D32").Value = 5.94
$ws.Range("E32").Value = 2931203187
$ws.Range("F32").Value = 234133852
$ws.Range("G32").Value = -1.03097

$ws.Range("B33").Value = "ETC"
$ws.Range("C33").Value = "Ethereum Classic"
$ws.Range("D33").Value = 19.71
$ws.Range("E33").Value = 2824819646
$ws.Range("F33").Value = 83148308
$ws.Range("G33").Value = -0.3235

$ws.Range("B34").Value = "APT"
$ws.Range("C34").Value = "Aptos"
$ws.Range("D34").Value = 8.970000000000001
$ws.Range("E34").Value = 2767436232
$ws.Range("F34").Value = 163140287
$ws.Range("G34").Value = -2.09505

$ws.Range("B35").Value = "XMR"
$ws.Range("C35").Value = "Monero"
$ws.Range("D35").Value = 152.48
$ws.Range("E35").Value = 2766137249
$ws.Range("F35").Value = 70156116
$ws.Range("G35").Value = -0.32068

$ws.Range("B36").Value = "LDO"
$ws.Range("C36").Value = "Lido DAO"
$ws.Range("D36").Value = 3.06
$ws.Range("E36").Value = 2717197183
$ws.Range("F36").Value = 173217378
$ws.Range("G36").Value = -6.53557

$ws.Range("B37").Value = "HBAR"
$ws.Range("C37").Value = "Hedera"
$ws.Range("D37").Value = 0.079376
$ws.Range("E37").Value = 2669707565
$ws.Range("F37").Value = 69649264
$ws.Range("G37").Value = -0.8775500000000001

$ws.Range("B38").Value = "IMX"
$ws.Range("C38").Value = "Immutable"
$ws.Range("D38").Value = 1.96
$ws.Range("E38").Value = 2592097856
$ws.Range("F38").Value = 58094130
$ws.Range("G38").Value = 0.52605

$ws.Range("B39").Value = "KAS"
$ws.Range("C39").Value = "Kaspa"
$ws.Range("D39").Value = 0.110256
$ws.Range("E39").Value = 2459912405
$ws.Range("F39").Value = 19045231
$ws.Range("G39").Value = 1.86428

$ws.Range("B40").Value = "STX"
$ws.Range("C40").Value = "Stacks"
$ws.Range("D40").Value = 1.66
$ws.Range("E40").Value = 2367965630
$ws.Range("F40").Value = 107741983
$ws.Range("G40").Value = 12.30051

$ws.Range("B41").Value = "CRO"
$ws.Range("C41").Value = "Cronos"
$ws.Range("D41").Value = 0.088865
$ws.Range("E41").Value = 2351065126
$ws.Range("F41").Value = 13427901
$ws.Range("G41").Value = -0.71765

$ws.Range("B42").Value = "ARB"
$ws.Range("C42").Value = "Arbitrum"
$ws.Range("D42").Value = 1.83
$ws.Range("E42").Value = 2316251069
$ws.Range("F42").Value = 669478108
$ws.Range("G42").Value = 2.70562

$ws.Range("B43").Value = "TUSD"
$ws.Range("C43").Value = "TrueUSD"
$ws.Range("D43").Value = 1.001
$ws.Range("E43").Value = 2200597722
$ws.Range("F43").Value = 137261838
$ws.Range("G43").Value = 0.08049000000000001

$ws.Range("B44").Value = "TIA"
$ws.Range("C44").Value = "Celestia"
$ws.Range("D44").Value = 15.09
$ws.Range("E44").Value = 2189609758
$ws.Range("F44").Value = 239569105
$ws.Range("G44").Value = 0.67301

$ws.Range("B45").Value = "VET"
$ws.Range("C45").Value = "VeChain"
$ws.Range("D45").Value = 0.02981264
$ws.Range("E45").Value = 2169111215
$ws.Range("F45").Value = 47585818
$ws.Range("G45").Value = 0.0332

$ws.Range("B46").Value = "MNT"
$ws.Range("C46").Value = "Mantle"
$ws.Range("D46").Value = 0.5990799999999999
$ws.Range("E46").Value = 1876671709
$ws.Range("F46").Value = 51696366
$ws.Range("G46").Value = 1.19627

$ws.Range("B47").Value = "FDUSD"
$ws.Range("C47").Value = "First Digital USD"
$ws.Range("D47").Value = 1.001
$ws.Range("E47").Value = 1841494564
$ws.Range("F47").Value = 459625980
$ws.Range("G47").Value = -0.01856

$ws.Range("B48").Value = "QNT"
$ws.Range("C48").Value = "Quant"
$ws.Range("D48").Value = 117.97
$ws.Range("E48").Value = 1715181715
$ws.Range("F48").Value = 25262778
$ws.Range("G48").Value = -0.18771

$ws.Range("B49").Value = "MKR"
$ws.Range("C49").Value = "Maker"
$ws.Range("D49").Value = 1792.66
$ws.Range("E49").Value = 1648380282
$ws.Range("F49").Value = 82960908
$ws.Range("G49").Value = 4.33082

$ws.Range("B50").Value = "GRT"
$ws.Range("C50").Value = "The Graph"
$ws.Range("D50").Value = 0.170412
$ws.Range("E50").Value = 1594815052
$ws.Range("F50").Value = 54716756
$ws.Range("G50").Value = 0.20431

$ws.Range("B51").Value = "SEI"
$ws.Range("C51").Value = "Sei"
$ws.Range("D51").Value = 0.667311
$ws.Range("E51").Value = 1535573064
$ws.Range("F51").Value = 565905294
$ws.Range("G51").Value = 1.9943
